$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rebuild the MASUK / PULANG header block (rows 10-12) so it has a new
# "cek waktu" indicator column inside each of the MASUK and PULANG
# groups, and the trailing JAM KERJA / SHIFT columns move two columns
# to the right (J,K -> L,M) to make room.
# ---------------------------------------------------------------------

# Break the old merges we are about to rebuild.
$ws.Range("D10:F10").UnMerge()
$ws.Range("G10:I10").UnMerge()
$ws.Range("J10:J11").UnMerge()
$ws.Range("K10:K11").UnMerge()

# --- Row 10 header labels ---
$ws.Range("H10").Value2 = $ws.Range("G10").Value2   # PULANG -> H10
$ws.Range("G10").ClearContents()
$ws.Range("L10").Value2 = $ws.Range("J10").Value2   # JAM KERJA -> L10
$ws.Range("M10").Value2 = $ws.Range("K10").Value2   # SHIFT -> M10
$ws.Range("J10").ClearContents()
$ws.Range("K10").ClearContents()

# --- Row 11 sub-headers (Jam / Verifikasi / Mesin, x2) ---
$ws.Range("K11").Value2 = $ws.Range("I11").Value2   # Mesin -> K11
$ws.Range("J11").Value2 = $ws.Range("H11").Value2   # Verifikasi -> J11
$ws.Range("H11").Value2 = $ws.Range("G11").Value2   # Jam -> H11
$ws.Range("G11").Value2 = $ws.Range("F11").Value2   # Mesin -> G11
$ws.Range("F11").Value2 = $ws.Range("E11").Value2   # Verifikasi -> F11
$ws.Range("E11").ClearContents()                    # D11 keeps "Jam"
$ws.Range("I11").ClearContents()                    # H11 keeps "Jam"

# --- Row 12 placeholder tokens ---
$ws.Range("M12").Value2 = "[jadwal_nama]"
$ws.Range("L12").Value2 = "[jam_kerja]"
$ws.Range("K12").Value2 = "[nama_mesin_pulang]"
$ws.Range("J12").Value2 = $ws.Range("H12").Value2
$ws.Range("I12").Value2 = "[cek_waktu_pulang]"
$ws.Range("H12").Value2 = $ws.Range("G12").Value2
$ws.Range("G12").Value2 = "[nama_mesin_masuk]"
$ws.Range("F12").Value2 = $ws.Range("E12").Value2
$ws.Range("E12").Value2 = "[cek_waktu_masuk]"

# --- Re-merge the header band with its new extents ---
$ws.Range("D10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:L11").Merge()
$ws.Range("M10:M11").Merge()
$ws.Range("D11:E11").Merge()
$ws.Range("H11:I11").Merge()

# --- Number format + alignment for the new "cek waktu" / "jam kerja" cells ---
$ws.Range("E10:E12,I10:I12,L10:L12").HorizontalAlignment = -4108
$ws.Range("E10:E12,I10:I12,L10:L12").VerticalAlignment = -4108
$ws.Range("E12,I12,L12").NumberFormat = "0.0"

# --- Column widths: the two new "cek waktu" columns are narrow, and the
#     shifted Mesin/Jam-kerja columns pick up new widths. ---
$ws.Columns("E:E").ColumnWidth = 4.6640625
$ws.Columns("I:I").ColumnWidth = 4.6640625
$ws.Columns("K:K").ColumnWidth = 19.21875
$ws.Columns("L:L").ColumnWidth = 11.44140625

# ---------------------------------------------------------------------
# Conditional formatting: highlight a "masuk" check value greater than 0
# and a "pulang" check value less than 0 with the classic red
# light-red-fill / dark-red-text warning style.
# ---------------------------------------------------------------------
$ws.Range("E10:E12").FormatConditions.Delete()
$cond1 = $ws.Range("E10:E12").FormatConditions.Add(1, 3, "0")
$cond1.Interior.Color = 13551615
$cond1.Font.Color = 10233776

$ws.Range("I10:I12").FormatConditions.Delete()
$cond2 = $ws.Range("I10:I12").FormatConditions.Add(1, 6, "0")
$cond2.Interior.Color = 13551615
$cond2.Font.Color = 10233776

# Selection ends up on L7 in the saved file.
$ws.Range("L7").Select()
